# Weekly update: insert a new daily-price record as the new row 371 in the
# "Hortaliza, Feria Lagunitas de Puerto Montt - Repollo" sheet, pushing the
# existing rows 371:432 down to 372:433.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 371 (copies formatting from the
# row being pushed down, same as Excel's native "Insert Row" behaviour).
$ws.Rows("371:371").Insert()

# Fill the newly inserted row with the new weekly record. Most fields repeat
# the values of the (now shifted) row below it; only the date, volume,
# min/max/weighted-avg price, origin and $/Kg price are new for this entry.
$ws.Range("A371").Value = 4
$ws.Range("B371").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C371").Value = "Los Lagos"
$ws.Range("D371").Value = 44694
$ws.Range("E371").Value = 10
$ws.Range("F371").Value = 100112006
$ws.Range("G371").Value = "Repollo"
$ws.Range("H371").Value = "Crespo record"
$ws.Range("I371").Value = "Primera"
$ws.Range("J371").Value = 1000
$ws.Range("K371").Value = 1800
$ws.Range("L371").Value = 2000
$ws.Range("M371").Value = 1900
$ws.Range("N371").Value = "$/unidad"
$ws.Range("O371").Value = "Región Metropolitana"
$ws.Range("P371").Value = 1900
$ws.Range("Q371").Value = 1
$ws.Range("R371").Value = "Hortaliza"
